$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New legend entries for "kiwoom real sw 2nd day"
$ws.Range("L7").Value = "방과후(2학기)"

$ws.Range("L8").Value = "웹툰"
$ws.Range("M8").Value = "화(15:40~17:10)"

$ws.Range("L9").Value = "생명과학"
$ws.Range("M9").Value = "금(15:00~16:00) or 16:10~17:10)"
$ws.Range("M9").WrapText = $true

# Column M width for readability of the new wrapped text
$ws.Columns.Item(13).ColumnWidth = 35.625

# Update last active selection cell to match the newly edited area
$ws.Range("K15").Select() | Out-Null
